# Sprint2.xlsx - "Upload excel sheets for working hours"
#
# Summary of the edit:
#   - B1 (hours worked) changes from 123 to 121; the dependent formulas in
#     B2 (=SUM(B1*50)) and B4 (=SUM(B2:B3)) recalculate automatically.
#   - A6 used to hold a plain-text link to the old TUNTIKIRJAUKSET wiki
#     page; it now points at the new "Resurssit-ja-tyotunnit" wiki page
#     and is turned into a real, styled (built-in "Hyperlink") hyperlink.
#   - The selected cell moves from B4 to D9.
#   - Column A widens to fit the longer URL text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hours worked for the week drives the two dependent SUM formulas ---
$ws.Range("B1").Value = 121

# --- Turn A6 into a hyperlink pointing at the new wiki page ---
$newUrl = "https://github.com/DigiaMinions/Project/wiki/Resurssit-ja-ty%C3%B6tunnit"
$ws.Range("A6").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("A6"), $newUrl)

# --- Column A needs to widen now that it holds the longer URL text ---
$ws.Columns("A").ColumnWidth = 69.67

# --- Move the live selection to D9 ---
$ws.Range("D9").Select()
